# Actualización automática 2025-09-12 16:50:09
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "VENTAS POR GRUPO"
# ---------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("E10").Value = 633.12
$wsGrupo.Range("H10").Value = 1037.85

$wsGrupo.Range("M13").Value = 5802.23

$wsGrupo.Range("H41").Value = 1010.7
$wsGrupo.Range("I41").Value = 26.1

$wsGrupo.Range("E60").Value = "2 de 58"
$wsGrupo.Range("H60").Value = "3 de 58"
$wsGrupo.Range("I60").Value = "2 de 58"
$wsGrupo.Range("M60").Value = "5 de 58"

# ---------------------------------------------------------------
# Sheet "VENTA MENSUAL"
# ---------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F10").Value = 1670.97
$wsMensual.Range("F13").Value = 5802.23
$wsMensual.Range("F41").Value = 1036.8
$wsMensual.Range("F60").Value = 21885.76

# Column D widened from 13 to 14 characters.
# (ColumnWidth applies a +0.8333.. offset when round-tripped through OOXML,
# so back that out to land exactly on the target stored width.)
$wsMensual.Columns.Item(4).ColumnWidth = 13.166666666666666

# ---------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumplimiento.Range("D4").Value = 1972.24
$wsCumplimiento.Range("E4").Value = -668.2113934184001
$wsCumplimiento.Range("F4").Value = 1.512420808903924

$wsCumplimiento.Range("D6").Value = 3236.55
$wsCumplimiento.Range("E6").Value = -2386.70580317996
$wsCumplimiento.Range("F6").Value = 3.808403954643183

$wsCumplimiento.Range("D7").Value = 585.78
$wsCumplimiento.Range("E7").Value = 123.588813030059
$wsCumplimiento.Range("F7").Value = 0.8257763651856203

$wsCumplimiento.Range("D12").Value = 10912.95
$wsCumplimiento.Range("E12").Value = 21491.85
$wsCumplimiento.Range("F12").Value = 0.3367695526588654

$wsCumplimiento.Range("D15").Value = 21885.76
$wsCumplimiento.Range("E15").Value = 28598.0070510252
$wsCumplimiento.Range("F15").Value = 0.4335207390106114
